# The commit swaps the contents of ppt/theme/theme1.xml (the "Integral" /
# "Red Violet" theme used by the slide master) and ppt/theme/theme2.xml
# (the "Office Theme" used by the notes master) -- theme1.xml ends up
# holding the Office Theme color scheme and theme2.xml ends up holding the
# Integral / Red Violet color scheme. The font scheme (fontScheme) and the
# format scheme (fmtScheme) are byte-for-byte identical between the two
# theme parts already, so the only observable difference after the swap is
# the color scheme applied to the deck's (single) design/theme, which is
# reachable through Slides/Designs -> SlideMaster -> Theme.

function ToBGR($rrggbb) {
    $r = ($rrggbb -band 0xFF0000) -shr 16
    $g = ($rrggbb -band 0x00FF00) -shr 8
    $b = ($rrggbb -band 0x0000FF)
    return ($b -shl 16) -bor ($g -shl 8) -bor $r
}

$p  = $ppt.ActivePresentation
$d  = $p.Designs.Item(1)
$sm = $d.SlideMaster
$cs = $sm.Theme.ThemeColorScheme

# Target palette = the former "Office Theme" color scheme (theme2.xml),
# now applied to theme1.xml. ThemeColorScheme order is:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$cs.Item(1).RGB  = ToBGR 0x000000   # dk1
$cs.Item(2).RGB  = ToBGR 0xFFFFFF   # lt1
$cs.Item(3).RGB  = ToBGR 0x44546A   # dk2
$cs.Item(4).RGB  = ToBGR 0xE7E6E6   # lt2
$cs.Item(5).RGB  = ToBGR 0x5B9BD5   # accent1
$cs.Item(6).RGB  = ToBGR 0xED7D31   # accent2
$cs.Item(7).RGB  = ToBGR 0xA5A5A5   # accent3
$cs.Item(8).RGB  = ToBGR 0xFFC000   # accent4
$cs.Item(9).RGB  = ToBGR 0x4472C4   # accent5
$cs.Item(10).RGB = ToBGR 0x70AD47   # accent6
$cs.Item(11).RGB = ToBGR 0x0563C1   # hlink
$cs.Item(12).RGB = ToBGR 0x954F72   # folHlink
